$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$numRows = 24

# Columns B:F (5 columns)
$arrBF = New-Object "object[,]" $numRows,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.040601056377843
$arrBF[0,2] = 1.044622031454461
$arrBF[0,3] = 1.03891922087721
$arrBF[0,4] = 1.053383609371474
$arrBF[1,0] = 1.02
$arrBF[1,1] = 1.042116107210576
$arrBF[1,2] = 1.0457687725996
$arrBF[1,3] = 1.040225825274041
$arrBF[1,4] = 1.054768806110744
$arrBF[2,0] = 1.02
$arrBF[2,1] = 1.043095084045744
$arrBF[2,2] = 1.04650941127792
$arrBF[2,3] = 1.04107031227016
$arrBF[2,4] = 1.055664136147774
$arrBF[3,0] = 1.02
$arrBF[3,1] = 1.043506327272124
$arrBF[3,2] = 1.046820450120802
$arrBF[3,3] = 1.041425107081118
$arrBF[3,4] = 1.056040302738354
$arrBF[4,0] = 1.02
$arrBF[4,1] = 1.043575358352176
$arrBF[4,2] = 1.046872655947639
$arrBF[4,3] = 1.041484665508734
$arrBF[4,4] = 1.05610344943045
$arrBF[5,0] = 1.02
$arrBF[5,1] = 1.043100580340588
$arrBF[5,2] = 1.046513568668099
$arrBF[5,3] = 1.041075053942981
$arrBF[5,4] = 1.055669163401579
$arrBF[6,0] = 1.02
$arrBF[6,1] = 1.041113361522263
$arrBF[6,2] = 1.045009866412437
$arrBF[6,3] = 1.039360998833051
$arrBF[6,4] = 1.053851950416018
$arrBF[7,0] = 1.02
$arrBF[7,1] = 1.037600856255405
$arrBF[7,2] = 1.042349378090465
$arrBF[7,3] = 1.03633291137544
$arrBF[7,4] = 1.05064199289539
$arrBF[8,0] = 1.02
$arrBF[8,1] = 1.035251469531577
$arrBF[8,2] = 1.040568183670831
$arrBF[8,3] = 1.034308657859507
$arrBF[8,4] = 1.048496423206938
$arrBF[9,0] = 1.02
$arrBF[9,1] = 1.034232222795397
$arrBF[9,2] = 1.039795053622449
$arrBF[9,3] = 1.033430744045676
$arrBF[9,4] = 1.047565960016393
$arrBF[10,0] = 1.02
$arrBF[10,1] = 1.033853327476458
$arrBF[10,2] = 1.039507593316746
$arrBF[10,3] = 1.033104430792118
$arrBF[10,4] = 1.047220124628235
$arrBF[11,0] = 1.02
$arrBF[11,1] = 1.033934615609805
$arrBF[11,2] = 1.039569267535452
$arrBF[11,3] = 1.033174436021593
$arrBF[11,4] = 1.047294317558069
$arrBF[12,0] = 1.02
$arrBF[12,1] = 1.034200909398271
$arrBF[12,2] = 1.039771297919174
$arrBF[12,3] = 1.033403775356293
$arrBF[12,4] = 1.047537377682459
$arrBF[13,0] = 1.02
$arrBF[13,1] = 1.034364941535171
$arrBF[13,2] = 1.039895737517041
$arrBF[13,3] = 1.033545049932919
$arrBF[13,4] = 1.047687105694179
$arrBF[14,0] = 1.02
$arrBF[14,1] = 1.035319071779018
$arrBF[14,2] = 1.04061945406294
$arrBF[14,3] = 1.034366892056942
$arrBF[14,4] = 1.048558144527166
$arrBF[15,0] = 1.02
$arrBF[15,1] = 1.03591704477842
$arrBF[15,2] = 1.041072920000211
$arrBF[15,3] = 1.034882032630621
$arrBF[15,4] = 1.04910413979111
$arrBF[16,0] = 1.02
$arrBF[16,1] = 1.036265645102559
$arrBF[16,2] = 1.041337239832302
$arrBF[16,3] = 1.035182370924411
$arrBF[16,4] = 1.049422473312345
$arrBF[17,0] = 1.02
$arrBF[17,1] = 1.036384477328026
$arrBF[17,2] = 1.041427335846411
$arrBF[17,3] = 1.035284755915514
$arrBF[17,4] = 1.049530993911629
$arrBF[18,0] = 1.02
$arrBF[18,1] = 1.035852907365505
$arrBF[18,2] = 1.0410242859717
$arrBF[18,3] = 1.034826776882556
$arrBF[18,4] = 1.049045573813182
$arrBF[19,0] = 1.02
$arrBF[19,1] = 1.034122500917502
$arrBF[19,2] = 1.039711812931177
$arrBF[19,3] = 1.033336246683081
$arrBF[19,4] = 1.047465808661493
$arrBF[20,0] = 1.02
$arrBF[20,1] = 1.033032775550566
$arrBF[20,2] = 1.0388849550116
$arrBF[20,3] = 1.032397832417051
$arrBF[20,4] = 1.046471271659126
$arrBF[21,0] = 1.02
$arrBF[21,1] = 1.033610628574498
$arrBF[21,2] = 1.039323446760952
$arrBF[21,3] = 1.032895425187939
$arrBF[21,4] = 1.046998617833584
$arrBF[22,0] = 1.02
$arrBF[22,1] = 1.035881888880332
$arrBF[22,2] = 1.041046262149683
$arrBF[22,3] = 1.034851744992481
$arrBF[22,4] = 1.049072037679065
$arrBF[23,0] = 1.02
$arrBF[23,1] = 1.038510248195776
$arrBF[23,2] = 1.043038484693539
$arrBF[23,3] = 1.037116693248241
$arrBF[23,4] = 1.051472804254958
$ws.Range("B2:F25").Value = $arrBF

# Columns I:N (6 columns)
$arrIN = New-Object "object[,]" $numRows,6
$arrIN[0,0] = 1.03500504795004
$arrIN[0,1] = 1.045686782448603
$arrIN[0,2] = 1.047392564274517
$arrIN[0,3] = 1.041705884398602
$arrIN[0,4] = 1.05612972849557
$arrIN[0,5] = 1.04717177762163
$arrIN[1,0] = 1.035307742326496
$arrIN[1,1] = 1.046844697082699
$arrIN[1,2] = 1.048349705849354
$arrIN[1,3] = 1.042821291055424
$arrIN[1,4] = 1.057326494006294
$arrIN[1,5] = 1.048331336627321
$arrIN[2,0] = 1.035501631538478
$arrIN[2,1] = 1.047592214494134
$arrIN[2,2] = 1.048967086084688
$arrIN[2,3] = 1.043541521930854
$arrIN[2,4] = 1.058099368677673
$arrIN[2,5] = 1.049079915599222
$arrIN[3,0] = 1.035582671512554
$arrIN[3,1] = 1.04790606186453
$arrIN[3,2] = 1.049226169169924
$arrIN[3,3] = 1.043843950107604
$arrIN[3,4] = 1.058423928161887
$arrIN[3,5] = 1.049394208668883
$arrIN[4,0] = 1.035596250904371
$arrIN[4,1] = 1.047958734394451
$arrIN[4,2] = 1.049269643301048
$arrIN[4,3] = 1.043894708347562
$arrIN[4,4] = 1.058478402338797
$arrIN[4,5] = 1.049446955999838
$arrIN[5,0] = 1.035502716247522
$arrIN[5,1] = 1.047596409738285
$arrIN[5,2] = 1.04897054978174
$arrIN[5,3] = 1.043545564387034
$arrIN[5,4] = 1.058103706853968
$arrIN[5,5] = 1.049084116801102
$arrIN[6,0] = 1.035107754687422
$arrIN[6,1] = 1.0460784667529
$arrIN[6,2] = 1.047716442020009
$arrIN[6,3] = 1.042083157379873
$arrIN[6,4] = 1.05653449728485
$arrIN[6,5] = 1.047564018162566
$arrIN[7,0] = 1.034396589088534
$arrIN[7,1] = 1.043390158428069
$arrIN[7,2] = 1.045491381552203
$arrIN[7,3] = 1.039494403535749
$arrIN[7,4] = 1.053757534928438
$arrIN[7,5] = 1.044871892131561
$arrIN[8,0] = 1.033912161603633
$arrIN[8,1] = 1.041588534886198
$arrIN[8,2] = 1.043997539670896
$arrIN[8,3] = 1.037760310304046
$arrIN[8,4] = 1.051897959741803
$arrIN[8,5] = 1.043067710077612
$arrIN[9,0] = 1.033699928898111
$arrIN[9,1] = 1.040806101384582
$arrIN[9,2] = 1.043348145064591
$arrIN[9,3] = 1.037007398466262
$arrIN[9,4] = 1.051090709629191
$arrIN[9,5] = 1.042284165430676
$arrIN[10,0] = 1.033620722708178
$arrIN[10,1] = 1.040515115692788
$arrIN[10,2] = 1.043106542297149
$arrIN[10,3] = 1.03672742102637
$arrIN[10,4] = 1.050790547880423
$arrIN[10,5] = 1.041992766505824
$arrIN[11,0] = 1.033637729640366
$arrIN[11,1] = 1.040577549286497
$arrIN[11,2] = 1.043158384587499
$arrIN[11,3] = 1.036787491378297
$arrIN[11,4] = 1.050854947866176
$arrIN[11,5] = 1.042055288762397
$arrIN[12,0] = 1.033693389321298
$arrIN[12,1] = 1.040782055698443
$arrIN[12,2] = 1.043328182065183
$arrIN[12,3] = 1.036984261854275
$arrIN[12,4] = 1.051065904589803
$arrIN[12,5] = 1.042260085596905
$arrIN[13,0] = 1.03372763353156
$arrIN[13,1] = 1.040908011599149
$arrIN[13,2] = 1.043432748226694
$arrIN[13,3] = 1.037105457046382
$arrIN[13,4] = 1.05119584030791
$arrIN[13,5] = 1.042386220369437
$arrIN[14,0] = 1.033926194522944
$arrIN[14,1] = 1.041640412999064
$arrIN[14,2] = 1.044040583640692
$arrIN[14,3] = 1.037810235110195
$arrIN[14,4] = 1.051951490713486
$arrIN[14,5] = 1.043119661863348
$arrIN[15,0] = 1.034050083202982
$arrIN[15,1] = 1.04209920356983
$arrIN[15,2] = 1.044421175902414
$arrIN[15,3] = 1.038251773721387
$arrIN[15,4] = 1.052424939131325
$arrIN[15,5] = 1.043579103969348
$arrIN[16,0] = 1.034122106963115
$arrIN[16,1] = 1.042366585090764
$arrIN[16,2] = 1.044642923142489
$arrIN[16,3] = 1.038509119509837
$arrIN[16,4] = 1.052700897013522
$arrIN[16,5] = 1.043846865202711
$arrIN[17,0] = 1.034146624827394
$arrIN[17,1] = 1.042457717613182
$arrIN[17,2] = 1.044718491693714
$arrIN[17,3] = 1.038596834680779
$arrIN[17,4] = 1.052794958369605
$arrIN[17,5] = 1.043938127143765
$arrIN[18,0] = 1.034036815795011
$arrIN[18,1] = 1.042050002824662
$arrIN[18,2] = 1.044380367416739
$arrIN[18,3] = 1.038204421147806
$arrIN[18,4] = 1.052374162975028
$arrIN[18,5] = 1.043529833353479
$arrIN[19,0] = 1.033677009259996
$arrIN[19,1] = 1.040721843512581
$arrIN[19,2] = 1.04327819171365
$arrIN[19,3] = 1.036926326529495
$arrIN[19,4] = 1.051003791787457
$arrIN[19,5] = 1.042199787902833
$arrIN[20,0] = 1.0334486225423
$arrIN[20,1] = 1.039884719460801
$arrIN[20,2] = 1.042582958516568
$arrIN[20,3] = 1.036120926777071
$arrIN[20,4] = 1.050140370263197
$arrIN[20,5] = 1.04136147503888
$arrIN[21,0] = 1.033569900294681
$arrIN[21,1] = 1.040328692045854
$arrIN[21,2] = 1.042951730028777
$arrIN[21,3] = 1.036548058186045
$arrIN[21,4] = 1.05059826050444
$arrIN[21,5] = 1.041806078115928
$arrIN[22,0] = 1.034042811502363
$arrIN[22,1] = 1.042072235213663
$arrIN[22,2] = 1.044398807774526
$arrIN[22,3] = 1.038225818344056
$arrIN[22,4] = 1.052397107142083
$arrIN[22,5] = 1.043552097315022
$arrIN[23,0] = 1.034582253344939
$arrIN[23,1] = 1.044086785315725
$arrIN[23,2] = 1.046068437899152
$arrIN[23,3] = 1.04016509028063
$arrIN[23,4] = 1.054476879500962
$arrIN[23,5] = 1.045569508309302
$ws.Range("I2:N25").Value = $arrIN
